$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 46000
$ws.Range("D8").Value = 162.23
$ws.Range("E8").Value = 160.99
$ws.Range("F8").Value = 170.99
$ws.Range("G8").Value = 161.11

$ws.Range("A9").Value = 46000
$ws.Range("D9").Value = 162.23
$ws.Range("E9").Value = 160.99
$ws.Range("F9").Value = 170.99
$ws.Range("G9").Value = 161.11

$ws.Range("A10").Value = 46000
$ws.Range("D10").Value = 164.17
$ws.Range("E10").Value = 163.32
$ws.Range("F10").Value = 173.32
$ws.Range("G10").Value = 163.84

$ws.Range("A11").Value = 45997
$ws.Range("D11").Value = 162.59
$ws.Range("E11").Value = 161.45
$ws.Range("F11").Value = 171.45
$ws.Range("G11").Value = 161.57

$ws.Range("A12").Value = 45997
$ws.Range("D12").Value = 162.59
$ws.Range("E12").Value = 161.45
$ws.Range("F12").Value = 171.45
$ws.Range("G12").Value = 161.57

$ws.Range("A13").Value = 45997
$ws.Range("D13").Value = 164.79
$ws.Range("E13").Value = 163.95
$ws.Range("F13").Value = 173.95
$ws.Range("G13").Value = 164.47

$ws.Range("A17").Value = 46000
$ws.Range("D17").Value = 167.08
$ws.Range("E17").Value = 166.69
$ws.Range("F17").Value = 176.69

$ws.Range("A18").Value = 45997
$ws.Range("D18").Value = 167.71
$ws.Range("E18").Value = 167.02
$ws.Range("F18").Value = 177.02

$ws.Range("A22").Value = 46000
$ws.Range("D22").Value = 163.26
$ws.Range("E22").Value = 162.5
$ws.Range("F22").Value = 172.1
$ws.Range("G22").Value = 163.66

$ws.Range("A23").Value = 46000
$ws.Range("D23").Value = 169.19
$ws.Range("E23").Value = 167.54
$ws.Range("F23").Value = 177.54
$ws.Range("G23").Value = "N/A"

$ws.Range("A24").Value = 46000
$ws.Range("D24").Value = 168.95
$ws.Range("E24").Value = 167.87
$ws.Range("F24").Value = 177.87
$ws.Range("G24").Value = "N/A"

$ws.Range("A25").Value = 46000
$ws.Range("D25").Value = 169.56
$ws.Range("E25").Value = 167.3
$ws.Range("F25").Value = 177.3
$ws.Range("G25").Value = 167.07

$ws.Range("A26").Value = 46000
$ws.Range("D26").Value = 168.36
$ws.Range("E26").Value = 168.75
$ws.Range("F26").Value = 178.75
$ws.Range("G26").Value = "N/A"

$ws.Range("A27").Value = 45997
$ws.Range("D27").Value = 163.84
$ws.Range("E27").Value = 163.24
$ws.Range("F27").Value = 172.84
$ws.Range("G27").Value = 164.4

$ws.Range("A28").Value = 45997
$ws.Range("D28").Value = 169.82
$ws.Range("E28").Value = 168.17
$ws.Range("F28").Value = 178.17
$ws.Range("G28").Value = "N/A"

$ws.Range("A29").Value = 45997
$ws.Range("D29").Value = 169.58
$ws.Range("E29").Value = 168.49
$ws.Range("F29").Value = 178.49
$ws.Range("G29").Value = "N/A"

$ws.Range("A30").Value = 45997
$ws.Range("D30").Value = 170.19
$ws.Range("E30").Value = 167.93
$ws.Range("F30").Value = 177.93
$ws.Range("G30").Value = 167.7

$ws.Range("A31").Value = 45997
$ws.Range("D31").Value = 168.99
$ws.Range("E31").Value = 169.39
$ws.Range("F31").Value = 179.39
$ws.Range("G31").Value = "N/A"

$ws.Range("A35").Value = 46000
$ws.Range("D35").Value = 162.44
$ws.Range("E35").Value = 160.33
$ws.Range("F35").Value = 169.33

$ws.Range("A36").Value = 45997
$ws.Range("D36").Value = 163.06
$ws.Range("E36").Value = 160.95
$ws.Range("F36").Value = 169.95

$ws.Range("A40").Value = 46000
$ws.Range("D40").Value = 168.47
$ws.Range("E40").Value = 167.37
$ws.Range("F40").Value = 177.37

$ws.Range("A41").Value = 46000
$ws.Range("D41").Value = 168.18
$ws.Range("E41").Value = 167.79
$ws.Range("F41").Value = 177.79

$ws.Range("A42").Value = 45997
$ws.Range("D42").Value = 169.09
$ws.Range("E42").Value = 167.68
$ws.Range("F42").Value = 177.68

$ws.Range("A43").Value = 45997
$ws.Range("D43").Value = 168.8
$ws.Range("E43").Value = 168.09
$ws.Range("F43").Value = 178.09

$ws.Range("A47").Value = 46000
$ws.Range("D47").Value = 163.54
$ws.Range("E47").Value = 162.56
$ws.Range("F47").Value = 172.56

$ws.Range("A48").Value = 46000
$ws.Range("D48").Value = 163.35
$ws.Range("E48").Value = 162.65
$ws.Range("F48").Value = 172.65

$ws.Range("A49").Value = 45997
$ws.Range("D49").Value = 163.65
$ws.Range("E49").Value = 162.9
$ws.Range("F49").Value = 172.9

$ws.Range("A50").Value = 45997
$ws.Range("D50").Value = 163.48
$ws.Range("E50").Value = 162.99
$ws.Range("F50").Value = 172.99

$ws.Range("A54").Value = 46000
$ws.Range("D54").Value = 177.99
$ws.Range("E54").Value = 178.19
$ws.Range("F54").Value = 188.19

$ws.Range("A55").Value = 46000
$ws.Range("D55").Value = 166.14
$ws.Range("E55").Value = 173.05
$ws.Range("F55").Value = 183.05

$ws.Range("A56").Value = 46000
$ws.Range("D56").Value = 168.62
$ws.Range("E56").Value = "N/A"
$ws.Range("F56").Value = "N/A"

$ws.Range("A57").Value = 46000
$ws.Range("D57").Value = 167.57
$ws.Range("E57").Value = 167.32
$ws.Range("F57").Value = "N/A"

$ws.Range("A58").Value = 46000
$ws.Range("D58").Value = 163.47
$ws.Range("E58").Value = 163.37
$ws.Range("F58").Value = 173.37

$ws.Range("A59").Value = 46000
$ws.Range("D59").Value = 169.98
$ws.Range("E59").Value = 175.41
$ws.Range("F59").Value = "N/A"

$ws.Range("A60").Value = 45997
$ws.Range("D60").Value = 178.64
$ws.Range("E60").Value = 178.77
$ws.Range("F60").Value = 188.77

$ws.Range("A61").Value = 45997
$ws.Range("D61").Value = 166.76
$ws.Range("E61").Value = 173.69
$ws.Range("F61").Value = 183.69

$ws.Range("A62").Value = 45997
$ws.Range("D62").Value = 169.24
$ws.Range("E62").Value = "N/A"
$ws.Range("F62").Value = "N/A"

$ws.Range("A63").Value = 45997
$ws.Range("D63").Value = 168.2
$ws.Range("E63").Value = 167.96
$ws.Range("F63").Value = "N/A"

$ws.Range("A64").Value = 45997
$ws.Range("D64").Value = 164.1
$ws.Range("E64").Value = 164.01
$ws.Range("F64").Value = 174.01

$ws.Range("A65").Value = 45997
$ws.Range("D65").Value = 170.61
$ws.Range("E65").Value = 176.02
$ws.Range("F65").Value = "N/A"
